$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg)
$rows = @(
  @{ Row = 2;  D = 44209; M = 58;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 },
  @{ Row = 3;  D = 44589; M = 60;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 },
  @{ Row = 4;  D = 44587; M = 165; N = 6500; O = 7000; P = 6742; R = "Provincia de Linares"; S = 3371 },
  @{ Row = 5;  D = 44606; M = 45;  N = 7000; O = 7000; P = 7000; R = "Provincia de Linares"; S = 3500 },
  @{ Row = 6;  D = 44214; M = 48;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 },
  @{ Row = 7;  D = 44592; M = 30;  N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 },
  @{ Row = 8;  D = 44211; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 },
  @{ Row = 9;  D = 44628; M = 40;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 },
  @{ Row = 10; D = 44588; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 },
  @{ Row = 11; D = 44582; M = 150; N = 6000; O = 6500; P = 6233; R = "Provincia de Curicó"; S = 3116 },
  @{ Row = 12; D = 44614; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 },
  @{ Row = 13; D = 44585; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 },
  @{ Row = 14; D = 44627; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 },
  @{ Row = 15; D = 44586; M = 80;  N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó"; S = 3500 }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 4).Value = $r.D    # D: Fecha
  $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Volumen
  $ws.Cells.Item($r.Row, 14).Value = $r.N   # N: Precio minimo
  $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Precio maximo
  $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio promedio ponderado
  $ws.Cells.Item($r.Row, 18).Value = $r.R   # R: Origen
  $ws.Cells.Item($r.Row, 19).Value = $r.S   # S: Precio $/Kg
}
